$wb = $excel.ActiveWorkbook

# Software sheet: rename the VEP software identifier (SMAHT_SOFTWARE_VEP -> SMAHT_SOFTWARE_VEPX).
$wsSoftware = $wb.Worksheets.Item("Software")
$wsSoftware.Range("A2").Value = "SMAHT_SOFTWARE_VEPX"
$wsSoftware.Range("A2").Select() | Out-Null

# Workflow sheet: update references to the renamed software identifier.
$wsWorkflow = $wb.Worksheets.Item("Workflow")
$wsWorkflow.Range("D2").Value = "SMAHT_SOFTWARE_VEPX"
$wsWorkflow.Range("D3").Value = "SMAHT_SOFTWARE_VEPX|SMAHT_SOFTWARE_FASTQC"
$wsWorkflow.Range("D3").Select() | Out-Null
$wsWorkflow.Activate() | Out-Null
